# Primeira versao do capitulo do modelo de estimacao dos juros da divida tecnica
#
# The only content edit captured by this script is on slide 1, shape
# "CaixaDeTexto 6": the word "instância" becomes "instâncias" and, as a
# side effect of that mid-sentence retype, PowerPoint splits the single
# run that used to hold the whole sentence into three runs:
#   "Cada uma " | "das instâncias " | "utiliza as métricas ... estimado."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

$oldFragment = "das instância "
$newFragment = "das instâncias "

$fragStart = $para.Text.IndexOf($oldFragment)
$seg = $para.Characters($fragStart + 1, $oldFragment.Length)
$seg.Text = $newFragment
